# Apply the edits described by the diff between Testdata/TC_26.xlsx revisions.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# 1. Rename the worksheet/tab from "My Series" to "Data"
#    (workbook.xml: <sheet name="My Series" .../> -> <sheet name="Data" .../>)
$ws.Name = "Data"

# 2. Number format for the data column (B27:B36) changes from
#    "0.00000000" to "###0.00000000"
$ws.Range("B27:B36").NumberFormat = "###0.00000000"

# 3. A11 label text changes from "Function Description" to "Function Information"
$ws.Range("A11").Value = "Function Information"

# 4. B21 (Kurtosis) value changes slightly
$ws.Range("B21").Value = 0.2499825759175085
